$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.265
$ws.Range("C7").Value = -12.675
$ws.Range("A10").Value = -21.736
$ws.Range("A12").Value = -21.696
$ws.Range("C15").Value = -13.636
$ws.Range("A18").Value = -22.035
$ws.Range("E18").Value = 16.412
$ws.Range("E19").Value = 16.538
$ws.Range("C20").Value = -12.183
$ws.Range("E27").Value = 16.357
$ws.Range("C29").Value = -11.987
$ws.Range("C30").Value = -13.347
$ws.Range("C31").Value = -13.519
$ws.Range("A37").Value = -20.029
$ws.Range("C40").Value = -12.782
$ws.Range("E42").Value = 16.576
$ws.Range("E44").Value = 16.542
$ws.Range("E47").Value = 16.399
$ws.Range("A55").Value = -21.868
$ws.Range("E58").Value = 16.541
$ws.Range("A68").Value = -21.736
$ws.Range("C68").Value = -11.001
$ws.Range("E73").Value = 16.546
$ws.Range("C76").Value = -12.72
$ws.Range("A77").Value = -20.843
$ws.Range("A78").Value = -20.134
$ws.Range("C87").Value = -13.197
$ws.Range("C88").Value = -13.089
$ws.Range("E95").Value = 17.399
$ws.Range("C96").Value = -12.586
$ws.Range("C98").Value = -13.201
$ws.Range("C101").Value = -12.747
$ws.Range("E101").Value = 16.701
$ws.Range("C102").Value = -13.091
